$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product titles
$ws.Range("A2").Value = "PigaphoneProXL"
$ws.Range("A3").Value = "PigBook Air"

# Update prices
$ws.Range("B2").Value = 80000
$ws.Range("B3").Value = 140000
